# Adds Czech translation rows for the new "clone build" feature on the
# "Import" sheet of the translations workbook (rows 540-553), mirroring
# the existing "vape" rows but for "build".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

$newRows = @(
    @("lab.build.preview.description", "Popis"),
    @("lab.build.preview.atomizer", "Atomizér"),
    @("lab.build.preview.coil", "Spirálka"),
    @("lab.build.preview.cotton", "Vata"),
    @("lab.build.preview.ohm", "Odpor"),
    @("lab.build.preview.coilOffset", "Pozice spirálky"),
    @("lab.build.preview.cottonOffset", "Množství vaty"),
    @("lab.build.preview.coils", "Počet spirálek"),
    @("lab.build.button.clone", "Klonovat"),
    @("lab.build.button.index", "Detail buildu"),
    @("lab.build.preview", "Náhled buildu"),
    @("lab.build.preview.preview.title", "Náhled buildu"),
    @("lab.build.preview.preview.subtitle", "Zde vidíte veškeré dostupné informace o buildu."),
    @("lab.build.clone.title", "Klon buildu")
)

$lastRow = 539
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i

    # Copy formatting (style) from the row above so the new rows keep the
    # same look (style index 1: wrapped text, 10pt font) as the rest of
    # the table.
    $ws.Range("A" + ($r - 1) + ":C" + ($r - 1)).Copy()
    $ws.Range("A" + $r + ":C" + $r).PasteSpecial(-4122)

    $ws.Range("A" + $r).Value2 = "cs"
    $ws.Range("B" + $r).Value2 = $newRows[$i][0]
    $ws.Range("C" + $r).Value2 = $newRows[$i][1]
}

$ws.Application.CutCopyMode = $false

$lastDataRow = $startRow + $newRows.Count - 1
$selRow = $lastDataRow - 5

[void]$ws.Range("B" + $selRow).Select()

Write-Host ("Added rows {0}-{1}" -f $startRow, $lastDataRow)
